$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -2
$ws.Range("F6").Value = -4
$ws.Range("F7").Value = -4
$ws.Range("F9").Value = -6
$ws.Range("F11").Value = -2
